$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 100.72727
$ws.Cells.Item(2, 9).Value = 104
$ws.Cells.Item(2, 10).Value = 95
$ws.Cells.Item(2, 11).Value = 104
$ws.Cells.Item(2, 12).Value = 95
$ws.Cells.Item(2, 13).Value = 9
$ws.Cells.Item(2, 14).Value = -321
$ws.Cells.Item(29, 8).Value = 2700
$ws.Cells.Item(29, 10).Value = 2700
$ws.Cells.Item(29, 12).Value = 8100
$ws.Cells.Item(29, 14).Value = -8662
$ws.Cells.Item(33, 8).Value = 9618205
$ws.Cells.Item(33, 9).Value = 16667809
$ws.Cells.Item(33, 11).Value = 16667809
$ws.Cells.Item(33, 13).Value = -16667580
$ws.Cells.Item(43, 8).Value = 13201.75
$ws.Cells.Item(43, 10).Value = 13201.75
$ws.Cells.Item(43, 12).Value = 13201.75
$ws.Cells.Item(43, 14).Value = -13339.75
$ws.Cells.Item(88, 8).Value = 3371.2307
$ws.Cells.Item(88, 9).Value = 5218.5
$ws.Cells.Item(88, 11).Value = 5218.5
$ws.Cells.Item(88, 13).Value = -4812.5
$ws.Cells.Item(91, 8).Value = 3371.2307
$ws.Cells.Item(91, 9).Value = 5218.5
$ws.Cells.Item(91, 11).Value = 5218.5
$ws.Cells.Item(91, 13).Value = -3814.5
$ws.Cells.Item(94, 8).Value = 951.4167
$ws.Cells.Item(94, 9).Value = 1173.8572
$ws.Cells.Item(94, 10).Value = 640
$ws.Cells.Item(94, 11).Value = 1173.8572
$ws.Cells.Item(94, 12).Value = 640
$ws.Cells.Item(94, 13).Value = -722.8571999999999
$ws.Cells.Item(94, 14).Value = -1542
$ws.Cells.Item(116, 8).Value = 5639.375
$ws.Cells.Item(116, 9).Value = 4324.9165
$ws.Cells.Item(116, 11).Value = 4324.9165
$ws.Cells.Item(116, 13).Value = -882.9165000000003
$ws.Cells.Item(132, 8).Value = 21187.867
$ws.Cells.Item(132, 9).Value = 18368.818
$ws.Cells.Item(132, 11).Value = 55106.454
$ws.Cells.Item(132, 13).Value = -52576.454
$ws.Cells.Item(137, 8).Value = 15344.066
$ws.Cells.Item(137, 9).Value = 1378.1818
$ws.Cells.Item(137, 11).Value = 4134.5454
$ws.Cells.Item(137, 13).Value = -1584.5454
$ws.Cells.Item(138, 8).Value = 2223.7415
$ws.Cells.Item(138, 9).Value = 2442.0454
$ws.Cells.Item(138, 11).Value = 7326.1362
$ws.Cells.Item(138, 13).Value = -2186.1362

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1994219.6
$ws.Cells.Item(5, 9).Value = 3081922.8
$ws.Cells.Item(5, 10).Value = 97.166664
$ws.Cells.Item(5, 11).Value = 3081922.8
$ws.Cells.Item(5, 12).Value = 97.166664
$ws.Cells.Item(5, 13).Value = -3081810.8
$ws.Cells.Item(5, 14).Value = -321.166664
$ws.Cells.Item(97, 8).Value = 851.5454999999999
$ws.Cells.Item(97, 9).Value = 851.5454999999999
$ws.Cells.Item(97, 11).Value = 851.5454999999999
$ws.Cells.Item(97, 13).Value = -355.5454999999999
$ws.Cells.Item(122, 8).Value = 1187048.1
$ws.Cells.Item(122, 9).Value = 1596017.8
$ws.Cells.Item(122, 11).Value = 4788053.4
$ws.Cells.Item(122, 13).Value = -4785603.4
$ws.Cells.Item(132, 8).Value = 2231117.2
$ws.Cells.Item(132, 9).Value = 2216.2896
$ws.Cells.Item(132, 11).Value = 6648.8688
$ws.Cells.Item(132, 13).Value = -4118.8688

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1994219.6
$ws.Cells.Item(4, 9).Value = 3081922.8
$ws.Cells.Item(4, 10).Value = 97.166664
$ws.Cells.Item(4, 11).Value = 3081922.8
$ws.Cells.Item(4, 12).Value = 97.166664
$ws.Cells.Item(4, 13).Value = -3081807.8
$ws.Cells.Item(4, 14).Value = -327.166664

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1821.5385
$ws.Cells.Item(22, 9).Value = 1200
$ws.Cells.Item(22, 11).Value = 1200
$ws.Cells.Item(22, 13).Value = -850
$ws.Cells.Item(31, 8).Value = 15808
$ws.Cells.Item(31, 9).Value = 909.2143
$ws.Cells.Item(31, 10).Value = 45605.57
$ws.Cells.Item(31, 11).Value = 909.2143
$ws.Cells.Item(31, 12).Value = 45605.57
$ws.Cells.Item(31, 13).Value = -614.2143
$ws.Cells.Item(31, 14).Value = -46195.57
$ws.Cells.Item(34, 8).Value = 15808
$ws.Cells.Item(34, 9).Value = 909.2143
$ws.Cells.Item(34, 10).Value = 45605.57
$ws.Cells.Item(34, 11).Value = 909.2143
$ws.Cells.Item(34, 12).Value = 45605.57
$ws.Cells.Item(34, 13).Value = -707.2143
$ws.Cells.Item(34, 14).Value = -46009.57
$ws.Cells.Item(58, 8).Value = 14071.525
$ws.Cells.Item(58, 9).Value = 5633.1665
$ws.Cells.Item(58, 11).Value = 5633.1665
$ws.Cells.Item(58, 13).Value = -5430.1665
$ws.Cells.Item(86, 8).Value = 11040.353
$ws.Cells.Item(86, 10).Value = 5880.5
$ws.Cells.Item(86, 12).Value = 5880.5
$ws.Cells.Item(86, 14).Value = -8126.5
$ws.Cells.Item(89, 8).Value = 11040.353
$ws.Cells.Item(89, 10).Value = 5880.5
$ws.Cells.Item(89, 12).Value = 29402.5
$ws.Cells.Item(89, 14).Value = -40634.5
$ws.Cells.Item(93, 8).Value = 12914.5
$ws.Cells.Item(93, 9).Value = 12914.5
$ws.Cells.Item(93, 11).Value = 12914.5
$ws.Cells.Item(93, 13).Value = -11042.5
$ws.Cells.Item(103, 8).Value = 10000
$ws.Cells.Item(103, 9).Value = 10000
$ws.Cells.Item(103, 11).Value = 10000
$ws.Cells.Item(103, 13).Value = -8828
$ws.Cells.Item(105, 8).Value = 5455.9565
$ws.Cells.Item(105, 9).Value = 8950.416999999999
$ws.Cells.Item(105, 10).Value = 1643.8182
$ws.Cells.Item(105, 11).Value = 8950.416999999999
$ws.Cells.Item(105, 12).Value = 1643.8182
$ws.Cells.Item(105, 13).Value = -7203.416999999999
$ws.Cells.Item(105, 14).Value = -5137.8182
$ws.Cells.Item(122, 8).Value = 2471.182
$ws.Cells.Item(122, 9).Value = 1454.7142
$ws.Cells.Item(122, 10).Value = 4250
$ws.Cells.Item(122, 11).Value = 4364.142599999999
$ws.Cells.Item(122, 12).Value = 12750
$ws.Cells.Item(122, 13).Value = -1914.142599999999
$ws.Cells.Item(122, 14).Value = -17650
$ws.Cells.Item(129, 8).Value = 94999
$ws.Cells.Item(129, 10).Value = 94999
$ws.Cells.Item(129, 12).Value = 94999
$ws.Cells.Item(129, 14).Value = -104999
$ws.Cells.Item(136, 8).Value = 14071.525
$ws.Cells.Item(136, 9).Value = 5633.1665
$ws.Cells.Item(136, 11).Value = 16899.4995
$ws.Cells.Item(136, 13).Value = -14349.4995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 9999
$ws.Cells.Item(26, 9).Value = 9999
$ws.Cells.Item(26, 11).Value = 29997
$ws.Cells.Item(26, 13).Value = -29709
$ws.Cells.Item(40, 8).Value = 405.05264
$ws.Cells.Item(40, 9).Value = 416.5625
$ws.Cells.Item(40, 10).Value = 343.66666
$ws.Cells.Item(40, 11).Value = 1666.25
$ws.Cells.Item(40, 12).Value = 1374.66664
$ws.Cells.Item(40, 13).Value = -1597.25
$ws.Cells.Item(40, 14).Value = -1512.66664
$ws.Cells.Item(131, 8).Value = 1356.6938
$ws.Cells.Item(131, 10).Value = 1443.1023
$ws.Cells.Item(131, 12).Value = 4329.3069
$ws.Cells.Item(131, 14).Value = -14409.3069

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 2727.2222
$ws.Cells.Item(43, 9).Value = 1717.25
$ws.Cells.Item(43, 10).Value = 3535.2
$ws.Cells.Item(43, 11).Value = 1717.25
$ws.Cells.Item(43, 12).Value = 3535.2
$ws.Cells.Item(43, 13).Value = -1566.25
$ws.Cells.Item(43, 14).Value = -3837.2
$ws.Cells.Item(46, 8).Value = 783.3333
$ws.Cells.Item(46, 9).Value = 783.3333
$ws.Cells.Item(46, 11).Value = 783.3333
$ws.Cells.Item(46, 13).Value = -627.3333
$ws.Cells.Item(134, 8).Value = 74600
$ws.Cells.Item(134, 10).Value = 74600
$ws.Cells.Item(134, 12).Value = 223800
$ws.Cells.Item(134, 14).Value = -228870

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 100003384
$ws.Cells.Item(16, 9).Value = 166670060
$ws.Cells.Item(16, 11).Value = 166670060
$ws.Cells.Item(16, 13).Value = -166669890
$ws.Cells.Item(105, 8).Value = 35000
$ws.Cells.Item(105, 10).Value = 35000
$ws.Cells.Item(105, 12).Value = 35000
$ws.Cells.Item(105, 14).Value = -41988
$ws.Cells.Item(132, 8).Value = 940923.7
$ws.Cells.Item(132, 9).Value = 2965.4375
$ws.Cells.Item(132, 10).Value = 1496750.9
$ws.Cells.Item(132, 11).Value = 8896.3125
$ws.Cells.Item(132, 12).Value = 4490252.699999999
$ws.Cells.Item(132, 13).Value = -6366.3125
$ws.Cells.Item(132, 14).Value = -4495312.699999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 7027
$ws.Cells.Item(29, 10).Value = 7027
$ws.Cells.Item(29, 12).Value = 7027
$ws.Cells.Item(29, 14).Value = -7607
$ws.Cells.Item(107, 8).Value = 1826.6
$ws.Cells.Item(107, 9).Value = 2075
$ws.Cells.Item(107, 10).Value = 833
$ws.Cells.Item(107, 11).Value = 6225
$ws.Cells.Item(107, 12).Value = 2499
$ws.Cells.Item(107, 13).Value = -4305
$ws.Cells.Item(107, 14).Value = -6339
$ws.Cells.Item(122, 8).Value = 1487387.9
$ws.Cells.Item(122, 9).Value = 2781352.2
$ws.Cells.Item(122, 11).Value = 8344056.600000001
$ws.Cells.Item(122, 13).Value = -8341606.600000001
$ws.Cells.Item(132, 8).Value = 8445.833000000001
$ws.Cells.Item(132, 9).Value = 4053
$ws.Cells.Item(132, 10).Value = 13936.875
$ws.Cells.Item(132, 11).Value = 12159
$ws.Cells.Item(132, 12).Value = 41810.625
$ws.Cells.Item(132, 13).Value = -9629
$ws.Cells.Item(132, 14).Value = -46870.625
